$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to track username/password/scenario/description; the
# "scenario" and "description" columns (C, D) are no longer used, so
# remove them completely (values + formatting).
$ws.Range("C1:D4").Clear()

# Append a new row of test credentials.
$ws.Range("A5").Value = "Admin"
$ws.Range("B5").Value = "admin123"

# Match the formatting already used by the rest of the username/password
# columns so the new row looks consistent with the existing data.
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)  # xlPasteFormats
